$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (C1 changes, D1/E1 are new)
$ws.Range("C1").Value = "Frecuencia del primer armonico"
$ws.Range("D1").Value = "Frecuencia del segundo armonico"
$ws.Range("E1").Value = "Frecuencia tercer armonico"

# Copy header style (bold + border + centered) from C1 onto the new D1:E1 headers
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate columns C (updated), D and E (new) for rows 2-57
$ws.Cells.Item(2, 3).Value = 379.2410950216681
$ws.Cells.Item(2, 4).Value = 189.8319416552158
$ws.Cells.Item(2, 5).Value = 561.0400591903608
$ws.Cells.Item(3, 3).Value = 385.4226292692238
$ws.Cells.Item(3, 4).Value = 193.0929211982448
$ws.Cells.Item(3, 5).Value = 576.9891242129361
$ws.Cells.Item(4, 3).Value = 407.7669902912621
$ws.Cells.Item(4, 4).Value = 135.0397175639891
$ws.Cells.Item(4, 5).Value = 270.0794351279783
$ws.Cells.Item(5, 3).Value = 362.6655714139029
$ws.Cells.Item(5, 4).Value = 181.9488653865901
$ws.Cells.Item(5, 5).Value = 543.3822774412156
$ws.Cells.Item(6, 3).Value = 359.5818815331013
$ws.Cells.Item(6, 4).Value = 180.5652342237709
$ws.Cells.Item(6, 5).Value = 715.137437088657
$ws.Cells.Item(7, 3).Value = 385.3596970971812
$ws.Cells.Item(7, 4).Value = 579.8906184265884
$ws.Cells.Item(7, 5).Value = 753.5549011358853
$ws.Cells.Item(8, 3).Value = 373.3121525019851
$ws.Cells.Item(8, 4).Value = 189.0389197776017
$ws.Cells.Item(8, 5).Value = 763.3042096902304
$ws.Cells.Item(9, 3).Value = 458.1580177653113
$ws.Cells.Item(9, 4).Value = 305.5633473585785
$ws.Cells.Item(9, 5).Value = 911.4539504441327
$ws.Cells.Item(10, 3).Value = 403.0423759507421
$ws.Cells.Item(10, 4).Value = 204.853314016661
$ws.Cells.Item(10, 5).Value = 807.5335023542202
$ws.Cells.Item(11, 3).Value = 283.3403745951273
$ws.Cells.Item(11, 4).Value = 567.2440501337842
$ws.Cells.Item(11, 5).Value = 848.8945218983245
$ws.Cells.Item(12, 3).Value = 378.5632839224627
$ws.Cells.Item(12, 4).Value = 568.3010262257694
$ws.Cells.Item(12, 5).Value = 190.6499429874575
$ws.Cells.Item(13, 3).Value = 489.8537045776311
$ws.Cells.Item(13, 4).Value = 316.1868806040584
$ws.Cells.Item(13, 5).Value = 654.0821142048135
$ws.Cells.Item(14, 3).Value = 185.062506678064
$ws.Cells.Item(14, 4).Value = 368.8428250881507
$ws.Cells.Item(14, 5).Value = 554.3327278555398
$ws.Cells.Item(15, 3).Value = 165.4900673124284
$ws.Cells.Item(15, 4).Value = 330.3234280085371
$ws.Cells.Item(15, 5).Value = 495.8134953209656
$ws.Cells.Item(16, 3).Value = 364.4405130473242
$ws.Cells.Item(16, 4).Value = 547.2504791390238
$ws.Cells.Item(16, 5).Value = 182.8099660916996
$ws.Cells.Item(17, 3).Value = 398.5260376833758
$ws.Cells.Item(17, 4).Value = 198.8458596954742
$ws.Cells.Item(17, 5).Value = 597.6500034763258
$ws.Cells.Item(18, 3).Value = 499.1539763113369
$ws.Cells.Item(18, 4).Value = 332.4873096446699
$ws.Cells.Item(18, 5).Value = 668.3587140439931
$ws.Cells.Item(19, 3).Value = 492.3076923076924
$ws.Cells.Item(19, 4).Value = 327.9720279720277
$ws.Cells.Item(19, 5).Value = 164.3356643356638
$ws.Cells.Item(20, 3).Value = 358.2966226138033
$ws.Cells.Item(20, 4).Value = 178.7287602265574
$ws.Cells.Item(20, 5).Value = 537.8644850010487
$ws.Cells.Item(21, 3).Value = 490.451841803193
$ws.Cells.Item(21, 4).Value = 317.6458311593447
$ws.Cells.Item(21, 5).Value = 652.4053010539501
$ws.Cells.Item(22, 3).Value = 498.0174464710544
$ws.Cells.Item(22, 4).Value = 166.3079188852389
$ws.Cells.Item(22, 5).Value = 332.1626826781467
$ws.Cells.Item(23, 3).Value = 429.4117647058829
$ws.Cells.Item(23, 4).Value = 644.1176470588234
$ws.Cells.Item(23, 5).Value = 216.9117647058829
$ws.Cells.Item(24, 3).Value = 444.7230929989551
$ws.Cells.Item(24, 4).Value = 703.8662486938347
$ws.Cells.Item(24, 5).Value = 563.8453500522464
$ws.Cells.Item(25, 3).Value = 465.7066145960462
$ws.Cells.Item(25, 4).Value = 310.1456350174922
$ws.Cells.Item(25, 5).Value = 155.2355381986818
$ws.Cells.Item(26, 3).Value = 485.43932008499
$ws.Cells.Item(26, 4).Value = 808.8988876390449
$ws.Cells.Item(26, 5).Value = 642.9196350456195
$ws.Cells.Item(27, 3).Value = 467.0237184391735
$ws.Cells.Item(27, 4).Value = 310.328997704667
$ws.Cells.Item(27, 5).Value = 155.2665136444784
$ws.Cells.Item(28, 3).Value = 600.2917578409924
$ws.Cells.Item(28, 4).Value = 399.7082421590085
$ws.Cells.Item(28, 5).Value = 200.583515681984
$ws.Cells.Item(29, 3).Value = 522.0804710500488
$ws.Cells.Item(29, 4).Value = 695.845600261694
$ws.Cells.Item(29, 5).Value = 869.0873405299317
$ws.Cells.Item(30, 3).Value = 483.9780711914136
$ws.Cells.Item(30, 4).Value = 645.2011427688985
$ws.Cells.Item(30, 5).Value = 806.4242143463825
$ws.Cells.Item(31, 3).Value = 834.2735671653991
$ws.Cells.Item(31, 4).Value = 555.6373150192139
$ws.Cells.Item(31, 5).Value = 694.301365383043
$ws.Cells.Item(32, 3).Value = 472.2924035708274
$ws.Cells.Item(32, 4).Value = 625.9053393970016
$ws.Cells.Item(32, 5).Value = 318.005726798046
$ws.Cells.Item(33, 3).Value = 486.7854947756605
$ws.Cells.Item(33, 4).Value = 246.6707641876665
$ws.Cells.Item(33, 5).Value = 732.6367547633681
$ws.Cells.Item(34, 3).Value = 524.7756055576046
$ws.Cells.Item(34, 4).Value = 263.1255379318818
$ws.Cells.Item(34, 5).Value = 1054.961268904463
$ws.Cells.Item(35, 3).Value = 522.2967757259066
$ws.Cells.Item(35, 4).Value = 261.5046612433935
$ws.Cells.Item(35, 5).Value = 1044.831067038775
$ws.Cells.Item(36, 3).Value = 583.9327322459731
$ws.Cells.Item(36, 4).Value = 293.7898895755243
$ws.Cells.Item(36, 5).Value = 1162.597507851281
$ws.Cells.Item(37, 3).Value = 531.1999999999998
$ws.Cells.Item(37, 4).Value = 265.8909090909092
$ws.Cells.Item(37, 5).Value = 795.3454545454542
$ws.Cells.Item(38, 3).Value = 496.0990045735807
$ws.Cells.Item(38, 4).Value = 993.2741458165192
$ws.Cells.Item(38, 5).Value = 744.3278629719307
$ws.Cells.Item(39, 3).Value = 500.6159152500613
$ws.Cells.Item(39, 4).Value = 255.2352796255236
$ws.Cells.Item(39, 5).Value = 747.9674796747968
$ws.Cells.Item(40, 3).Value = 488.2090503505419
$ws.Cells.Item(40, 4).Value = 244.741873804971
$ws.Cells.Item(40, 5).Value = 739.3244104525174
$ws.Cells.Item(41, 3).Value = 507.4080239720324
$ws.Cells.Item(41, 4).Value = 255.0357915764935
$ws.Cells.Item(41, 5).Value = 1013.484268353587
$ws.Cells.Item(42, 3).Value = 398.2523008273683
$ws.Cells.Item(42, 4).Value = 183.3224876824397
$ws.Cells.Item(42, 5).Value = 552.1985683740822
$ws.Cells.Item(43, 3).Value = 447.4327628361857
$ws.Cells.Item(43, 4).Value = 223.3088834555829
$ws.Cells.Item(43, 5).Value = 670.7416462917686
$ws.Cells.Item(44, 3).Value = 465.6416692814555
$ws.Cells.Item(44, 4).Value = 232.1932852212112
$ws.Cells.Item(44, 5).Value = 925.0078443677439
$ws.Cells.Item(45, 3).Value = 528.3806343906508
$ws.Cells.Item(45, 4).Value = 263.7729549248752
$ws.Cells.Item(45, 5).Value = 792.153589315526
$ws.Cells.Item(46, 3).Value = 464.0142166057849
$ws.Cells.Item(46, 4).Value = 233.3892783098036
$ws.Cells.Item(46, 5).Value = 701.7474577944513
$ws.Cells.Item(47, 3).Value = 283.7912391040491
$ws.Cells.Item(47, 4).Value = 569.7892530067302
$ws.Cells.Item(47, 5).Value = 847.8428776343371
$ws.Cells.Item(48, 3).Value = 269.1388426160638
$ws.Cells.Item(48, 4).Value = 537.9879771130582
$ws.Cells.Item(48, 5).Value = 808.5753603244739
$ws.Cells.Item(49, 3).Value = 610.3121711680114
$ws.Cells.Item(49, 4).Value = 304.4545773412838
$ws.Cells.Item(49, 5).Value = 916.8712732374606
$ws.Cells.Item(50, 3).Value = 557.7342047930279
$ws.Cells.Item(50, 4).Value = 279.9564270152505
$ws.Cells.Item(50, 5).Value = 1117.102396514161
$ws.Cells.Item(51, 3).Value = 171.805230703806
$ws.Cells.Item(51, 4).Value = 515.4156921114181
$ws.Cells.Item(51, 5).Value = 345.3115032957685
$ws.Cells.Item(52, 3).Value = 440.4157224349456
$ws.Cells.Item(52, 4).Value = 293.5062905368441
$ws.Cells.Item(52, 5).Value = 587.3251543330471
$ws.Cells.Item(53, 3).Value = 447.8138222849084
$ws.Cells.Item(53, 4).Value = 223.5543018335684
$ws.Cells.Item(53, 5).Value = 894.9224259520452
$ws.Cells.Item(54, 3).Value = 401.2488992074295
$ws.Cells.Item(54, 4).Value = 198.2227203586581
$ws.Cells.Item(54, 5).Value = 601.0727723961254
$ws.Cells.Item(55, 3).Value = 408.0934042423128
$ws.Cells.Item(55, 4).Value = 204.1689589828229
$ws.Cells.Item(55, 5).Value = 818.3874319946208
$ws.Cells.Item(56, 3).Value = 450.7628294036067
$ws.Cells.Item(56, 4).Value = 902.2191400832176
$ws.Cells.Item(56, 5).Value = 224.6879334257974
$ws.Cells.Item(57, 3).Value = 485.6061852278335
$ws.Cells.Item(57, 4).Value = 243.4610955749299
$ws.Cells.Item(57, 5).Value = 971.212370455667

Write-Host "Updated header row 1 and columns C/D/E for data rows 2-57"
